# Quarterly financial update for ABBV: inserts two new quarterly columns
# (D, E) ahead of the existing quarter columns, shifting prior data right,
# and populates the two new columns with the latest reported figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank columns before the existing "D" quarter column. This
# shifts all existing quarter data (old D:K) right to F:M.
$ws.Columns("D:E").Insert()

# The newly inserted columns are blank/unstyled; clone the number
# formatting/styling from the (now shifted) first data column - which still
# holds the same per-row style the old column D used to have - onto the two
# new columns for every row in the used range.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns with their reported values.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 8305000
$ws.Range("E8").Value = 8236000
$ws.Range("D9").Value = 2005000
$ws.Range("E9").Value = 1781000
$ws.Range("D10").Value = 6300000
$ws.Range("E10").Value = 6455000
$ws.Range("D12").Value = 6495000
$ws.Range("E12").Value = 1268000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 300000
$ws.Range("E14").Value = 55000
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 10746000
$ws.Range("E17").Value = 5077000
$ws.Range("D18").Value = -2441000
$ws.Range("E18").Value = 3159000
$ws.Range("D20").Value = 448000
$ws.Range("E20").Value = -59000
$ws.Range("D21").Value = -1551000
$ws.Range("E21").Value = 3535000
$ws.Range("D22").Value = 380000
$ws.Range("E22").Value = 339000
$ws.Range("D23").Value = -2373000
$ws.Range("E23").Value = 2761000
$ws.Range("D24").Value = -633000
$ws.Range("E24").Value = 14000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -1740000
$ws.Range("E26").Value = 2747000
$ws.Range("D27").Value = -1736000
$ws.Range("E27").Value = 2735000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -86000
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -448000
$ws.Range("E32").Value = 59000
$ws.Range("D33").Value = -1822000
$ws.Range("E33").Value = 2735000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -1822000
$ws.Range("E35").Value = 2735000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 7289000
$ws.Range("E41").Value = 8015000
$ws.Range("D42").Value = 772000
$ws.Range("E42").Value = 770000
$ws.Range("D43").Value = 5384000
$ws.Range("E43").Value = 5780000
$ws.Range("D44").Value = 1605000
$ws.Range("E44").Value = 1786000
$ws.Range("D45").Value = 1895000
$ws.Range("E45").Value = 2114000
$ws.Range("D46").Value = 16945000
$ws.Range("E46").Value = 18465000
$ws.Range("D47").Value = 1420000
$ws.Range("E47").Value = 1463000
$ws.Range("D48").Value = 2883000
$ws.Range("E48").Value = 2950000
$ws.Range("D49").Value = 36896000
$ws.Range("E49").Value = 42343000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1208000
$ws.Range("E52").Value = 943000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 59352000
$ws.Range("E54").Value = 66164000
$ws.Range("D57").Value = 5572000
$ws.Range("E57").Value = 11343000
$ws.Range("D58").Value = 5308000
$ws.Range("E58").Value = 4021000
$ws.Range("D59").Value = 6359000
$ws.Range("E59").Value = 23000
$ws.Range("D60").Value = 17239000
$ws.Range("E60").Value = 15387000
$ws.Range("D61").Value = 35002000
$ws.Range("E61").Value = 36487000
$ws.Range("D62").Value = 15557000
$ws.Range("E62").Value = 17211000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 67798000
$ws.Range("E66").Value = 69085000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 3368000
$ws.Range("E72").Value = 6789000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -8446000
$ws.Range("E76").Value = -2921000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -1822000
$ws.Range("E81").Value = 2735000
$ws.Range("D83").Value = 442000
$ws.Range("E83").Value = 435000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 3392000
$ws.Range("E89").Value = 4524000
$ws.Range("D91").Value = -123000
$ws.Range("E91").Value = -282000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -283000
$ws.Range("E94").Value = -963000
$ws.Range("D96").Value = -1451000
$ws.Range("E96").Value = -1461000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -3825000
$ws.Range("E100").Value = 916000
$ws.Range("D101").Value = -10000
$ws.Range("E101").Value = -9000
$ws.Range("D102").Value = -726000
$ws.Range("E102").Value = 4468000

# Row 29 ("Discontinued Operations") keeps its "NA" placeholders elsewhere
# in the row; only explicit values changed above for D29/E29.

$ws.Range("A1").Select()
